$wb = $excel.ActiveWorkbook

$wsInstructions = $wb.Worksheets.Item("Instructions")
$wsSalesForecast = $wb.Worksheets.Item("Sales Forecast")

# --- Update "Instructions" sheet ---
# The sheet is protected; unprotect so the cell values can be edited, then
# re-protect it to match the original state.
$wsInstructions.Unprotect()

# Insert a new note about not renaming the sheet right after the "fill out" note,
# and update the wording of the numbered validation rules to use camelCase field
# names (and tweak the model_year rule wording).

$wsInstructions.Range("A3").Value = 'Please do not alter the name of the "Sales Forecast" sheet'
$wsInstructions.Range("A5").Value = 'Please do not alter any of the header names in the "Sales Forecast" sheet'
$wsInstructions.Range("A7").Value = 'Please note that all of the fields in the "Sales Forecast" sheet are mandatory '
$wsInstructions.Range("A8").ClearContents()
$wsInstructions.Range("A9").Value = "Please note that:"
$wsInstructions.Range("A10").Value = '(1) "modelYear" should be a 4 digit integer'
$wsInstructions.Range("A11").Value = '(2) "make" should be no more than 250 characters'
$wsInstructions.Range("A12").Value = '(3) "modelName" should be no more than 250 characters'
$wsInstructions.Range("A13").Value = '(4) "type" should be exactly one of: BEV, PHEV, FCEV, EREV'
$wsInstructions.Range("A14").Value = '(5) "range" should be a real number with no more than 2 decimal places'
$wsInstructions.Range("A15").Value = '(6) "zevClass" should be a single, uppercase letter'
$wsInstructions.Range("A16").Value = '(7) "interiorVolume" should be a real number with no more than 2 decimal places'
$wsInstructions.Range("A17").Value = '(8) "totalSales" should be an integer'

$wsInstructions.Protect("", $true, $true, $true)

# --- Update "Sales Forecast" header row to camelCase field names ---
$wsSalesForecast.Range("A1").Value = "modelYear"
$wsSalesForecast.Range("B1").Value = "make"
$wsSalesForecast.Range("C1").Value = "modelName"
$wsSalesForecast.Range("D1").Value = "type"
$wsSalesForecast.Range("E1").Value = "range"
$wsSalesForecast.Range("F1").Value = "zevClass"
$wsSalesForecast.Range("G1").Value = "interiorVolume"
$wsSalesForecast.Range("H1").Value = "totalSales"
